$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Cells.Item(14, 1).Value = 7794000006515
$ws.Cells.Item(14, 2).Value = "Aderezo"
$ws.Cells.Item(14, 3).Value = "a base de mostaza"
$ws.Cells.Item(14, 4).Value = "original"
$ws.Cells.Item(14, 5).Value = "Savora"
$ws.Cells.Item(14, 6).Value = 60
$ws.Cells.Item(14, 7).Value = "gr."
$ws.Cells.Item(14, 8).Value = "Sobre"
$ws.Cells.Item(14, 9).Value = "Aderezos"
$ws.Cells.Item(14, 10).Value = "Argentina"
$ws.Cells.Item(14, 11).Value = 24
$ws.Cells.Item(14, 12).Value = $false
$ws.Cells.Item(14, 13).Value = $true
$ws.Cells.Item(14, 14).Value = "C:\VentaSoft\Imágenes de artículos\7794000006515.png"
$ws.Cells.Item(14, 15).Value = $true

# Row 15
$ws.Cells.Item(15, 1).Value = 1111
$ws.Cells.Item(15, 2).Value = "Salame"
$ws.Cells.Item(15, 3).Value = "pelado"
$ws.Cells.Item(15, 4).Value = "milán"
$ws.Cells.Item(15, 5).Value = "Fela"
$ws.Cells.Item(15, 6).Value = 1000
$ws.Cells.Item(15, 7).Value = "gr."
$ws.Cells.Item(15, 8).Value = "Horma"
$ws.Cells.Item(15, 9).Value = "Fiambres"
$ws.Cells.Item(15, 10).Value = "Argentina"
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = $true
$ws.Cells.Item(15, 13).Value = $true
$ws.Cells.Item(15, 14).Value = "C:\VentaSoft\Imágenes de artículos\1111.png"
$ws.Cells.Item(15, 15).Value = $true

# Row 16
$ws.Cells.Item(16, 1).Value = 2222
$ws.Cells.Item(16, 2).Value = "Salame"
$ws.Cells.Item(16, 3).Value = "tipo"
$ws.Cells.Item(16, 4).Value = "crespón"
$ws.Cells.Item(16, 5).Value = "Paladini"
$ws.Cells.Item(16, 6).Value = 1000
$ws.Cells.Item(16, 7).Value = "gr."
$ws.Cells.Item(16, 8).Value = "Horma"
$ws.Cells.Item(16, 9).Value = "Fiambres"
$ws.Cells.Item(16, 10).Value = "Argentina"
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = $true
$ws.Cells.Item(16, 13).Value = $true
$ws.Cells.Item(16, 14).Value = "C:\VentaSoft\Imágenes de artículos\2222.png"
$ws.Cells.Item(16, 15).Value = $true

# Row 17
$ws.Cells.Item(17, 1).Value = 3333
$ws.Cells.Item(17, 2).Value = "Queso"
$ws.Cells.Item(17, 3).Value = "tipo"
$ws.Cells.Item(17, 4).Value = "cremoso"
$ws.Cells.Item(17, 5).Value = "La Paulina"
$ws.Cells.Item(17, 6).Value = 1000
$ws.Cells.Item(17, 7).Value = "gr."
$ws.Cells.Item(17, 8).Value = "Horma"
$ws.Cells.Item(17, 9).Value = "Quesos"
$ws.Cells.Item(17, 10).Value = "Argentina"
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = $true
$ws.Cells.Item(17, 13).Value = $true
$ws.Cells.Item(17, 14).Value = "C:\VentaSoft\Imágenes de artículos\3333.png"
$ws.Cells.Item(17, 15).Value = $true

# Row 18
$ws.Cells.Item(18, 1).Value = 77990112
$ws.Cells.Item(18, 2).Value = "Licor"
$ws.Cells.Item(18, 3).Value = "crema"
$ws.Cells.Item(18, 4).Value = "café al coñac"
$ws.Cells.Item(18, 5).Value = "Tres Plumas"
$ws.Cells.Item(18, 6).Value = 200
$ws.Cells.Item(18, 7).Value = "ml."
$ws.Cells.Item(18, 8).Value = "Petaca"
$ws.Cells.Item(18, 9).Value = "Licores"
$ws.Cells.Item(18, 10).Value = "Argentina"
$ws.Cells.Item(18, 11).Value = 12
$ws.Cells.Item(18, 12).Value = $false
$ws.Cells.Item(18, 13).Value = $true
$ws.Cells.Item(18, 14).Value = "C:\VentaSoft\Imágenes de artículos\77990112.png"
$ws.Cells.Item(18, 15).Value = $true

# Row 19
$ws.Cells.Item(19, 1).Value = 3086123206873
$ws.Cells.Item(19, 2).Value = "Encendedor"
$ws.Cells.Item(19, 3).Value = "no recargable"
$ws.Cells.Item(19, 4).Value = "colores varios"
$ws.Cells.Item(19, 5).Value = "BX7"
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = "und."
$ws.Cells.Item(19, 8).Value = "Suelto"
$ws.Cells.Item(19, 9).Value = "Encendedores"
$ws.Cells.Item(19, 10).Value = "Francia"
$ws.Cells.Item(19, 11).Value = 50
$ws.Cells.Item(19, 12).Value = $false
$ws.Cells.Item(19, 13).Value = $false
$ws.Cells.Item(19, 14).Value = "C:\VentaSoft\Imágenes de artículos\3086123206873.png"
$ws.Cells.Item(19, 15).Value = $true

# Row 20
$ws.Cells.Item(20, 1).Value = 7790036973036
$ws.Cells.Item(20, 2).Value = "Puré de tomates"
$ws.Cells.Item(20, 3).Value = "sin conservantes"
$ws.Cells.Item(20, 4).Value = "de la huerta"
$ws.Cells.Item(20, 5).Value = "Baggio"
$ws.Cells.Item(20, 6).Value = 210
$ws.Cells.Item(20, 7).Value = "gr."
$ws.Cells.Item(20, 8).Value = "Tetra Brik"
$ws.Cells.Item(20, 9).Value = "Puré de tomates"
$ws.Cells.Item(20, 10).Value = "Argentina"
$ws.Cells.Item(20, 11).Value = 12
$ws.Cells.Item(20, 12).Value = $false
$ws.Cells.Item(20, 13).Value = $true
$ws.Cells.Item(20, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790036973036.png"
$ws.Cells.Item(20, 15).Value = $true

# Row 21
$ws.Cells.Item(21, 1).Value = 7790070336385
$ws.Cells.Item(21, 2).Value = "Fideos"
$ws.Cells.Item(21, 3).Value = "secos"
$ws.Cells.Item(21, 4).Value = "spaghetti Nº 7"
$ws.Cells.Item(21, 5).Value = "Luchetti"
$ws.Cells.Item(21, 6).Value = 500
$ws.Cells.Item(21, 7).Value = "gr."
$ws.Cells.Item(21, 8).Value = "Bolsa"
$ws.Cells.Item(21, 9).Value = "Fideos secos"
$ws.Cells.Item(21, 10).Value = "Argentina"
$ws.Cells.Item(21, 11).Value = 12
$ws.Cells.Item(21, 12).Value = $false
$ws.Cells.Item(21, 13).Value = $true
$ws.Cells.Item(21, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790070336385.png"
$ws.Cells.Item(21, 15).Value = $true

# Row 22
$ws.Cells.Item(22, 1).Value = 7790150006153
$ws.Cells.Item(22, 2).Value = "Café"
$ws.Cells.Item(22, 3).Value = "torrado molido"
$ws.Cells.Item(22, 4).Value = "clásico"
$ws.Cells.Item(22, 5).Value = "La Virginia"
$ws.Cells.Item(22, 6).Value = 500
$ws.Cells.Item(22, 7).Value = "gr."
$ws.Cells.Item(22, 8).Value = "Bolsa"
$ws.Cells.Item(22, 9).Value = "Cafés"
$ws.Cells.Item(22, 10).Value = "Argentina"
$ws.Cells.Item(22, 11).Value = 12
$ws.Cells.Item(22, 12).Value = $false
$ws.Cells.Item(22, 13).Value = $true
$ws.Cells.Item(22, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790150006153.png"
$ws.Cells.Item(22, 15).Value = $true

# Row 23
$ws.Cells.Item(23, 1).Value = 7790150430392
$ws.Cells.Item(23, 2).Value = "Bicarbonato"
$ws.Cells.Item(23, 3).Value = "de sodio"
$ws.Cells.Item(23, 4).Value = "en bolsa"
$ws.Cells.Item(23, 5).Value = "Alicante"
$ws.Cells.Item(23, 6).Value = 50
$ws.Cells.Item(23, 7).Value = "gr."
$ws.Cells.Item(23, 8).Value = "Bolsa"
$ws.Cells.Item(23, 9).Value = "Especias"
$ws.Cells.Item(23, 10).Value = "Argentina"
$ws.Cells.Item(23, 11).Value = 12
$ws.Cells.Item(23, 12).Value = $false
$ws.Cells.Item(23, 13).Value = $true
$ws.Cells.Item(23, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790150430392.png"
$ws.Cells.Item(23, 15).Value = $true

# Row 24
$ws.Cells.Item(24, 1).Value = 7790520010445
$ws.Cells.Item(24, 2).Value = "Espirales insecticidas"
$ws.Cells.Item(24, 3).Value = "repelentes"
$ws.Cells.Item(24, 4).Value = "contra mosquitos"
$ws.Cells.Item(24, 5).Value = "Raid"
$ws.Cells.Item(24, 6).Value = 12
$ws.Cells.Item(24, 7).Value = "und."
$ws.Cells.Item(24, 8).Value = "Caja"
$ws.Cells.Item(24, 9).Value = "Insecticidas"
$ws.Cells.Item(24, 10).Value = "Indonesia"
$ws.Cells.Item(24, 11).Value = 12
$ws.Cells.Item(24, 12).Value = $false
$ws.Cells.Item(24, 13).Value = $true
$ws.Cells.Item(24, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790520010445.png"
$ws.Cells.Item(24, 15).Value = $true

# Row 25
$ws.Cells.Item(25, 1).Value = 7790520016461
$ws.Cells.Item(25, 2).Value = "Desodorante"
$ws.Cells.Item(25, 3).Value = "para ambientes en aerosol"
$ws.Cells.Item(25, 4).Value = "caricias de algodón"
$ws.Cells.Item(25, 5).Value = "Glade"
$ws.Cells.Item(25, 6).Value = 360
$ws.Cells.Item(25, 7).Value = "cm3."
$ws.Cells.Item(25, 8).Value = "Lata"
$ws.Cells.Item(25, 9).Value = "Desodorantes"
$ws.Cells.Item(25, 10).Value = "Argentina"
$ws.Cells.Item(25, 11).Value = 6
$ws.Cells.Item(25, 12).Value = $false
$ws.Cells.Item(25, 13).Value = $true
$ws.Cells.Item(25, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790520016461.png"
$ws.Cells.Item(25, 15).Value = $true

# Row 26
$ws.Cells.Item(26, 1).Value = 7790520995360
$ws.Cells.Item(26, 2).Value = "Tabletas insecticidas"
$ws.Cells.Item(26, 3).Value = "repelentes"
$ws.Cells.Item(26, 4).Value = "contra mosquitos"
$ws.Cells.Item(26, 5).Value = "Raid"
$ws.Cells.Item(26, 6).Value = 24
$ws.Cells.Item(26, 7).Value = "und."
$ws.Cells.Item(26, 8).Value = "Caja"
$ws.Cells.Item(26, 9).Value = "Insecticidas"
$ws.Cells.Item(26, 10).Value = "Argentina"
$ws.Cells.Item(26, 11).Value = 24
$ws.Cells.Item(26, 12).Value = $false
$ws.Cells.Item(26, 13).Value = $true
$ws.Cells.Item(26, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790520995360.png"
$ws.Cells.Item(26, 15).Value = $true

# Row 27
$ws.Cells.Item(27, 1).Value = 7790639003574
$ws.Cells.Item(27, 2).Value = "Gaseosa"
$ws.Cells.Item(27, 3).Value = "indian tonic"
$ws.Cells.Item(27, 4).Value = "classic"
$ws.Cells.Item(27, 5).Value = "Cunnington"
$ws.Cells.Item(27, 6).Value = 1.5
$ws.Cells.Item(27, 7).Value = "lt."
$ws.Cells.Item(27, 8).Value = "Botella"
$ws.Cells.Item(27, 9).Value = "Gaseosas"
$ws.Cells.Item(27, 10).Value = "Argentina"
$ws.Cells.Item(27, 11).Value = 6
$ws.Cells.Item(27, 12).Value = $false
$ws.Cells.Item(27, 13).Value = $true
$ws.Cells.Item(27, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790639003574.png"
$ws.Cells.Item(27, 15).Value = $true

# Row 28
$ws.Cells.Item(28, 1).Value = 7791070000382
$ws.Cells.Item(28, 2).Value = "Papel higiénico"
$ws.Cells.Item(28, 3).Value = "simple hoja"
$ws.Cells.Item(28, 4).Value = """soft"""
$ws.Cells.Item(28, 5).Value = "Ecco"
$ws.Cells.Item(28, 6).Value = 12
$ws.Cells.Item(28, 7).Value = "und."
$ws.Cells.Item(28, 8).Value = "Bolsa"
$ws.Cells.Item(28, 9).Value = "Papeles Higiénicos"
$ws.Cells.Item(28, 10).Value = "Argentina"
$ws.Cells.Item(28, 11).Value = 4
$ws.Cells.Item(28, 12).Value = $false
$ws.Cells.Item(28, 13).Value = $false
$ws.Cells.Item(28, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791070000382.png"
$ws.Cells.Item(28, 15).Value = $false

# Row 29
$ws.Cells.Item(29, 1).Value = 7791600087128
$ws.Cells.Item(29, 2).Value = "Desodorante"
$ws.Cells.Item(29, 3).Value = "para hombre"
$ws.Cells.Item(29, 4).Value = "US"
$ws.Cells.Item(29, 5).Value = "Colbert"
$ws.Cells.Item(29, 6).Value = 150
$ws.Cells.Item(29, 7).Value = "ml."
$ws.Cells.Item(29, 8).Value = "Lata"
$ws.Cells.Item(29, 9).Value = "Desodorantes"
$ws.Cells.Item(29, 10).Value = "Argentina"
$ws.Cells.Item(29, 11).Value = 12
$ws.Cells.Item(29, 12).Value = $false
$ws.Cells.Item(29, 13).Value = $true
$ws.Cells.Item(29, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791600087128.png"
$ws.Cells.Item(29, 15).Value = $true

# Row 30
$ws.Cells.Item(30, 1).Value = 7791664000156
$ws.Cells.Item(30, 2).Value = "Ravioles"
$ws.Cells.Item(30, 3).Value = "rellenos de"
$ws.Cells.Item(30, 4).Value = "carne y espinaca"
$ws.Cells.Item(30, 5).Value = "La Italiana"
$ws.Cells.Item(30, 6).Value = 500
$ws.Cells.Item(30, 7).Value = "gr."
$ws.Cells.Item(30, 8).Value = "Blister"
$ws.Cells.Item(30, 9).Value = "Pastas frescas"
$ws.Cells.Item(30, 10).Value = "Argentina"
$ws.Cells.Item(30, 11).Value = 6
$ws.Cells.Item(30, 12).Value = $false
$ws.Cells.Item(30, 13).Value = $true
$ws.Cells.Item(30, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791664000156.png"
$ws.Cells.Item(30, 15).Value = $true

# Row 31
$ws.Cells.Item(31, 1).Value = 7791664000453
$ws.Cells.Item(31, 2).Value = "Ñoquis"
$ws.Cells.Item(31, 3).Value = "frescos"
$ws.Cells.Item(31, 4).Value = "de papa"
$ws.Cells.Item(31, 5).Value = "La Italiana"
$ws.Cells.Item(31, 6).Value = 500
$ws.Cells.Item(31, 7).Value = "gr."
$ws.Cells.Item(31, 8).Value = "Blister"
$ws.Cells.Item(31, 9).Value = "Pastas frescas"
$ws.Cells.Item(31, 10).Value = "Argentina"
$ws.Cells.Item(31, 11).Value = 6
$ws.Cells.Item(31, 12).Value = $false
$ws.Cells.Item(31, 13).Value = $true
$ws.Cells.Item(31, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791664000453.png"
$ws.Cells.Item(31, 15).Value = $true

# Row 32
$ws.Cells.Item(32, 1).Value = 7791684000934
$ws.Cells.Item(32, 2).Value = "Ravioles"
$ws.Cells.Item(32, 3).Value = "rellenos de"
$ws.Cells.Item(32, 4).Value = "jamón y muzzarella"
$ws.Cells.Item(32, 5).Value = "La Italiana"
$ws.Cells.Item(32, 6).Value = 500
$ws.Cells.Item(32, 7).Value = "gr."
$ws.Cells.Item(32, 8).Value = "Blister"
$ws.Cells.Item(32, 9).Value = "Pastas frescas"
$ws.Cells.Item(32, 10).Value = "Argentina"
$ws.Cells.Item(32, 11).Value = 12
$ws.Cells.Item(32, 12).Value = $false
$ws.Cells.Item(32, 13).Value = $true
$ws.Cells.Item(32, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791684000934.png"
$ws.Cells.Item(32, 15).Value = $true

# Row 33
$ws.Cells.Item(33, 1).Value = 7793913000139
$ws.Cells.Item(33, 2).Value = "Queso"
$ws.Cells.Item(33, 3).Value = "rallado"
$ws.Cells.Item(33, 4).Value = "sin gluten/TACC"
$ws.Cells.Item(33, 5).Value = "Tregar"
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = "und."
$ws.Cells.Item(33, 8).Value = "Sobre"
$ws.Cells.Item(33, 9).Value = "Quesos"
$ws.Cells.Item(33, 10).Value = "Argentina"
$ws.Cells.Item(33, 11).Value = 12
$ws.Cells.Item(33, 12).Value = $false
$ws.Cells.Item(33, 13).Value = $true
$ws.Cells.Item(33, 14).Value = "C:\VentaSoft\Imágenes de artículos\7793913000139.png"
$ws.Cells.Item(33, 15).Value = $true

# Row 34
$ws.Cells.Item(34, 1).Value = 7795018002902
$ws.Cells.Item(34, 2).Value = "Leche"
$ws.Cells.Item(34, 3).Value = "en polvo"
$ws.Cells.Item(34, 4).Value = "entera"
$ws.Cells.Item(34, 5).Value = "Quelech"
$ws.Cells.Item(34, 6).Value = 800
$ws.Cells.Item(34, 7).Value = "gr."
$ws.Cells.Item(34, 8).Value = "Bolsa"
$ws.Cells.Item(34, 9).Value = "Papeles Higiénicos"
$ws.Cells.Item(34, 10).Value = "Argentina"
$ws.Cells.Item(34, 11).Value = 12
$ws.Cells.Item(34, 12).Value = $false
$ws.Cells.Item(34, 13).Value = $true
$ws.Cells.Item(34, 14).Value = "C:\VentaSoft\Imágenes de artículos\7795018002902.png"
$ws.Cells.Item(34, 15).Value = $true

# Row 35
$ws.Cells.Item(35, 1).Value = 7798100200491
$ws.Cells.Item(35, 2).Value = "Cigarrillos"
$ws.Cells.Item(35, 3).Value = "rubios mentolados"
$ws.Cells.Item(35, 4).Value = "en caja"
$ws.Cells.Item(35, 5).Value = "Milenio"
$ws.Cells.Item(35, 6).Value = 20
$ws.Cells.Item(35, 7).Value = "und."
$ws.Cells.Item(35, 8).Value = "Caja"
$ws.Cells.Item(35, 9).Value = "Cigarrillos"
$ws.Cells.Item(35, 10).Value = "Argentina"
$ws.Cells.Item(35, 11).Value = 10
$ws.Cells.Item(35, 12).Value = $false
$ws.Cells.Item(35, 13).Value = $false
$ws.Cells.Item(35, 14).Value = "C:\VentaSoft\Imágenes de artículos\7798100200491.png"
$ws.Cells.Item(35, 15).Value = $true

# Row 36
$ws.Cells.Item(36, 1).Value = 7798100200583
$ws.Cells.Item(36, 2).Value = "Cigarrillos"
$ws.Cells.Item(36, 3).Value = "rubios convertibles"
$ws.Cells.Item(36, 4).Value = "en caja"
$ws.Cells.Item(36, 5).Value = "Mill"
$ws.Cells.Item(36, 6).Value = 20
$ws.Cells.Item(36, 7).Value = "und."
$ws.Cells.Item(36, 8).Value = "Caja"
$ws.Cells.Item(36, 9).Value = "Cigarrillos"
$ws.Cells.Item(36, 10).Value = "Argentina"
$ws.Cells.Item(36, 11).Value = 10
$ws.Cells.Item(36, 12).Value = $false
$ws.Cells.Item(36, 13).Value = $false
$ws.Cells.Item(36, 14).Value = "C:\VentaSoft\Imágenes de artículos\7798100200583.png"
$ws.Cells.Item(36, 15).Value = $true

# Row 37
$ws.Cells.Item(37, 1).Value = 7791600087012
$ws.Cells.Item(37, 2).Value = "Eau de toilette"
$ws.Cells.Item(37, 3).Value = "para hombre"
$ws.Cells.Item(37, 4).Value = "US"
$ws.Cells.Item(37, 5).Value = "Colbert"
$ws.Cells.Item(37, 6).Value = 60
$ws.Cells.Item(37, 7).Value = "ml."
$ws.Cells.Item(37, 8).Value = "Caja"
$ws.Cells.Item(37, 9).Value = "Perfumes/Colonias"
$ws.Cells.Item(37, 10).Value = "Argentina"
$ws.Cells.Item(37, 11).Value = 6
$ws.Cells.Item(37, 12).Value = $false
$ws.Cells.Item(37, 13).Value = $true
$ws.Cells.Item(37, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791600087012.png"
$ws.Cells.Item(37, 15).Value = $true

# Row 38
$ws.Cells.Item(38, 1).Value = 7791520009743
$ws.Cells.Item(38, 2).Value = "Talco"
$ws.Cells.Item(38, 3).Value = "para hombre"
$ws.Cells.Item(38, 4).Value = "original"
$ws.Cells.Item(38, 5).Value = "Veritas"
$ws.Cells.Item(38, 6).Value = 180
$ws.Cells.Item(38, 7).Value = "gr."
$ws.Cells.Item(38, 8).Value = "Botella"
$ws.Cells.Item(38, 9).Value = "Desodorantes"
$ws.Cells.Item(38, 10).Value = "Argentina"
$ws.Cells.Item(38, 11).Value = 6
$ws.Cells.Item(38, 12).Value = $false
$ws.Cells.Item(38, 13).Value = $true
$ws.Cells.Item(38, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791520009743.png"
$ws.Cells.Item(38, 15).Value = $true

# Row 39
$ws.Cells.Item(39, 1).Value = 7794000006065
$ws.Cells.Item(39, 2).Value = "Mayonesa"
$ws.Cells.Item(39, 3).Value = "libre de gluten"
$ws.Cells.Item(39, 4).Value = "clásica"
$ws.Cells.Item(39, 5).Value = "Hellmann's"
$ws.Cells.Item(39, 6).Value = 237
$ws.Cells.Item(39, 7).Value = "gr."
$ws.Cells.Item(39, 8).Value = "Pouch"
$ws.Cells.Item(39, 9).Value = "Aderezos"
$ws.Cells.Item(39, 10).Value = "Argentina"
$ws.Cells.Item(39, 11).Value = 24
$ws.Cells.Item(39, 12).Value = $false
$ws.Cells.Item(39, 13).Value = $true
$ws.Cells.Item(39, 14).Value = "C:\VentaSoft\Imágenes de artículos\7794000006065.png"
$ws.Cells.Item(39, 15).Value = $true

# Row 40
$ws.Cells.Item(40, 1).Value = 7792540260138
$ws.Cells.Item(40, 2).Value = "Azúcar"
$ws.Cells.Item(40, 3).Value = "común tipo ""a"""
$ws.Cells.Item(40, 4).Value = "clásica"
$ws.Cells.Item(40, 5).Value = "Ledesma"
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(40, 7).Value = "kg."
$ws.Cells.Item(40, 8).Value = "Bolsa"
$ws.Cells.Item(40, 9).Value = "Azúcar"
$ws.Cells.Item(40, 10).Value = "Argentina"
$ws.Cells.Item(40, 11).Value = 10
$ws.Cells.Item(40, 12).Value = $false
$ws.Cells.Item(40, 13).Value = $true
$ws.Cells.Item(40, 14).Value = "C:\VentaSoft\Imágenes de artículos\7792540260138.png"
$ws.Cells.Item(40, 15).Value = $true

# Row 41
$ws.Cells.Item(41, 1).Value = 7790639003895
$ws.Cells.Item(41, 2).Value = "Gaseosa"
$ws.Cells.Item(41, 3).Value = "indian tonic"
$ws.Cells.Item(41, 4).Value = "suave"
$ws.Cells.Item(41, 5).Value = "Cunnington"
$ws.Cells.Item(41, 6).Value = 1.5
$ws.Cells.Item(41, 7).Value = "lt."
$ws.Cells.Item(41, 8).Value = "Botella"
$ws.Cells.Item(41, 9).Value = "Gaseosas"
$ws.Cells.Item(41, 10).Value = "Argentina"
$ws.Cells.Item(41, 11).Value = 6
$ws.Cells.Item(41, 12).Value = $false
$ws.Cells.Item(41, 13).Value = $true
$ws.Cells.Item(41, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790639003895.png"
$ws.Cells.Item(41, 15).Value = $true

# Row 42
$ws.Cells.Item(42, 1).Value = 7791600174767
$ws.Cells.Item(42, 2).Value = "Desodorante"
$ws.Cells.Item(42, 3).Value = "para hombre"
$ws.Cells.Item(42, 4).Value = "noir"
$ws.Cells.Item(42, 5).Value = "Colbert"
$ws.Cells.Item(42, 6).Value = 250
$ws.Cells.Item(42, 7).Value = "ml."
$ws.Cells.Item(42, 8).Value = "Lata"
$ws.Cells.Item(42, 9).Value = "Desodorantes"
$ws.Cells.Item(42, 10).Value = "Argentina"
$ws.Cells.Item(42, 11).Value = 6
$ws.Cells.Item(42, 12).Value = $false
$ws.Cells.Item(42, 13).Value = $true
$ws.Cells.Item(42, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791600174767.png"
$ws.Cells.Item(42, 15).Value = $true

# Row 43
$ws.Cells.Item(43, 1).Value = 7791600192372
$ws.Cells.Item(43, 2).Value = "Desodorante"
$ws.Cells.Item(43, 3).Value = "para hombre"
$ws.Cells.Item(43, 4).Value = "verde"
$ws.Cells.Item(43, 5).Value = "Colbert"
$ws.Cells.Item(43, 6).Value = 150
$ws.Cells.Item(43, 7).Value = "ml."
$ws.Cells.Item(43, 8).Value = "Lata"
$ws.Cells.Item(43, 9).Value = "Desodorantes"
$ws.Cells.Item(43, 10).Value = "Argentina"
$ws.Cells.Item(43, 11).Value = 6
$ws.Cells.Item(43, 12).Value = $false
$ws.Cells.Item(43, 13).Value = $true
$ws.Cells.Item(43, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791600192372.png"
$ws.Cells.Item(43, 15).Value = $true

# Row 44
$ws.Cells.Item(44, 1).Value = 7791600192488
$ws.Cells.Item(44, 2).Value = "Desodorante"
$ws.Cells.Item(44, 3).Value = "para hombre"
$ws.Cells.Item(44, 4).Value = "verde"
$ws.Cells.Item(44, 5).Value = "Colbert"
$ws.Cells.Item(44, 6).Value = 250
$ws.Cells.Item(44, 7).Value = "ml."
$ws.Cells.Item(44, 8).Value = "Lata"
$ws.Cells.Item(44, 9).Value = "Desodorantes"
$ws.Cells.Item(44, 10).Value = "Argentina"
$ws.Cells.Item(44, 11).Value = 6
$ws.Cells.Item(44, 12).Value = $false
$ws.Cells.Item(44, 13).Value = $true
$ws.Cells.Item(44, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791600192488.png"
$ws.Cells.Item(44, 15).Value = $true

# Row 45
$ws.Cells.Item(45, 1).Value = 2001300000798
$ws.Cells.Item(45, 2).Value = "Huevos"
$ws.Cells.Item(45, 3).Value = "frescos medianos"
$ws.Cells.Item(45, 4).Value = "de granja"
$ws.Cells.Item(45, 5).Value = "Cismondi"
$ws.Cells.Item(45, 6).Value = 6
$ws.Cells.Item(45, 7).Value = "und."
$ws.Cells.Item(45, 8).Value = "Maple"
$ws.Cells.Item(45, 9).Value = "Huevos"
$ws.Cells.Item(45, 10).Value = "Argentina"
$ws.Cells.Item(45, 11).Value = 1
$ws.Cells.Item(45, 12).Value = $false
$ws.Cells.Item(45, 13).Value = $true
$ws.Cells.Item(45, 14).Value = "C:\VentaSoft\Imágenes de artículos\2001300000798.png"
$ws.Cells.Item(45, 15).Value = $false

# Row 46
$ws.Cells.Item(46, 1).Value = 7791070000078
$ws.Cells.Item(46, 2).Value = "Papel higiénico"
$ws.Cells.Item(46, 3).Value = "simple hoja"
$ws.Cells.Item(46, 4).Value = """soft"""
$ws.Cells.Item(46, 5).Value = "Campanita"
$ws.Cells.Item(46, 6).Value = 4
$ws.Cells.Item(46, 7).Value = "und."
$ws.Cells.Item(46, 8).Value = "Bolsa"
$ws.Cells.Item(46, 9).Value = "Papeles Higiénicos"
$ws.Cells.Item(46, 10).Value = "Argentina"
$ws.Cells.Item(46, 11).Value = 10
$ws.Cells.Item(46, 12).Value = $false
$ws.Cells.Item(46, 13).Value = $false
$ws.Cells.Item(46, 14).Value = "C:\VentaSoft\Imágenes de artículos\7791070000078.png"
$ws.Cells.Item(46, 15).Value = $true

$ws.Range("A46").NumberFormat = "0"
